$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> column letter -> new value (Price=D, Volume(1h)=E, Hora=G)
$updates = @{
  2 = @{ "D"="305.07"; "E"="1.46%"; "G"="16" }
  3 = @{ "D"="35.60"; "E"="0.61%"; "G"="16" }
  4 = @{ "D"="5.073"; "E"="0.49%"; "G"="16" }
  5 = @{ "D"="0.08056"; "E"="1.34%"; "G"="16" }
  6 = @{ "D"="1.913"; "E"="1.25%"; "G"="16" }
  7 = @{ "D"="4.180"; "E"="3.71%"; "G"="16" }
  8 = @{ "D"="7.842"; "E"="0.83%"; "G"="16" }
  9 = @{ "D"="0.9331"; "E"="0.61%"; "G"="16" }
  10 = @{ "D"="0.1366"; "E"="-1.23%"; "G"="16" }
  11 = @{ "D"="0.1892"; "E"="-0.09%"; "G"="16" }
  12 = @{ "D"="0.09176"; "E"="0.54%"; "G"="16" }
  13 = @{ "D"="0.03504"; "E"="0.23%"; "G"="16" }
  14 = @{ "D"="0.09907"; "E"="-0.03%"; "G"="16" }
  15 = @{ "D"="0.001421"; "E"="2.58%"; "G"="16" }
  16 = @{ "D"="0.006727"; "E"="13.33%"; "G"="16" }
  17 = @{ "D"="3.616"; "E"="2.86%"; "G"="16" }
  18 = @{ "D"="2.972"; "E"="1.25%"; "G"="16" }
  19 = @{ "D"="0.3427"; "E"="0.53%"; "G"="16" }
  20 = @{ "D"="0.1340"; "E"="3.55%"; "G"="16" }
  21 = @{ "D"="5.194"; "E"="2.96%"; "G"="16" }
  22 = @{ "G"="16" }
  23 = @{ "D"="0.04418"; "E"="-1.89%"; "G"="16" }
  24 = @{ "D"="0.001241"; "E"="2.09%"; "G"="16" }
  25 = @{ "D"="0.004690"; "E"="-1.46%"; "G"="16" }
  26 = @{ "D"="0.0001303"; "E"="6.06%"; "G"="16" }
  27 = @{ "D"="0.0003148"; "E"="4.80%"; "G"="16" }
  28 = @{ "G"="16" }
  29 = @{ "G"="16" }
  30 = @{ "G"="16" }
  31 = @{ "G"="16" }
  32 = @{ "G"="16" }
  33 = @{ "G"="16" }
  34 = @{ "G"="16" }
  35 = @{ "G"="16" }
  36 = @{ "G"="16" }
  37 = @{ "G"="16" }
  38 = @{ "G"="16" }
  39 = @{ "D"="0.01975"; "E"="5.88%"; "G"="16" }
  40 = @{ "D"="0.05073"; "E"="7.43%"; "G"="16" }
  41 = @{ "D"="0.01129"; "E"="13.54%"; "G"="16" }
  42 = @{ "D"="0.007635"; "G"="16" }
  43 = @{ "D"="0.1365"; "E"="3.32%"; "G"="16" }
  44 = @{ "D"="0.002104"; "E"="-0.13%"; "G"="16" }
  45 = @{ "D"="0.01077"; "E"="15.60%"; "G"="16" }
  46 = @{ "D"="0.00006388"; "E"="1.57%"; "G"="16" }
  47 = @{ "D"="0.00000000755"; "E"="0.51%"; "G"="16" }
  48 = @{ "E"="0.82%"; "G"="16" }
  49 = @{ "D"="0.001610"; "E"="-2.97%"; "G"="16" }
  50 = @{ "D"="0.00002114"; "E"="0.51%"; "G"="16" }
  51 = @{ "D"="0.0002013"; "E"="0.51%"; "G"="16" }
}

foreach ($row in $updates.Keys) {
  $cols = $updates[$row]
  foreach ($col in $cols.Keys) {
    $cellRef = "$col$row"
    $newVal = $cols[$col]
    $ws.Range($cellRef).Value = "'$newVal"
    $ws.Range($cellRef).Style = "Normal"
  }
}
